$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E10").Value = 355
$ws.Range("E25").Value = 184
$ws.Range("F35").Value = 63
$ws.Range("H35").Value = 63
$ws.Range("E46").Value = 230
$ws.Range("E47").Value = 333
$ws.Range("E48").Value = 152
